$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts existing A:D -> B:E)
$ws.Columns("A:A").Insert()

# Insert a new row 3 (Files tab) below the existing data row
$ws.Rows("3:3").Insert()

# Column A - tab name labels
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "FilesTab"

# Column B - main queries
$ws.Range("B2").Value = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
    WHERE a.pubmed_id IN ['31504139', '31765263'] 
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@
$ws.Range("B3").Value = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.pubmed_id IN ['31504139', '31765263'] 
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# Column C - stat query (same text reused on both rows)
$ws.Range("C2").Value = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
    WHERE a.pubmed_id IN ['31504139', '31765263'] 
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@
$ws.Range("C3").Value = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
    WHERE a.pubmed_id IN ['31504139', '31765263'] 
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# Column D / E - carry the existing filenames onto the new Files row
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("E3").Value = $ws.Range("E2").Value()

# Formatting: wrap text on the big query cells + trailing styled cell
$ws.Range("B2:C3").WrapText = $true
$ws.Range("C4").WrapText = $true

# Row heights
$ws.Rows("2:2").RowHeight = 195
$ws.Rows("3:3").RowHeight = 409.5

# Column widths (approximate AutoFit sizing)
$ws.Columns("A:A").ColumnWidth = 8
$ws.Columns("B:C").ColumnWidth = 75
$ws.Columns("D:D").ColumnWidth = 70
$ws.Columns("E:E").ColumnWidth = 28

# Selection / view
$ws.Activate()
$ws.Range("C12:C13").Select()
